$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D10 value from "-" to "+"
$ws.Range("D10").Value = "+"

# Update active selection to D11
$ws.Range("D11").Select()
